# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D, preserving it as literal text even
# when the string would otherwise be auto-parsed by Excel as a number
# (the sheet stores prices/amounts as plain text, e.g. '1.00' or '6.21').
function Set-TextValue($cell, [string]$text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $text
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '57.417.53'
$ws.Cells.Item(2, 5).Value = '  +1.21%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '3.013.14'
$ws.Cells.Item(3, 5).Value = '  +0.16%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '508.79'
$ws.Cells.Item(5, 5).Value = '  -0.40%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '139.48'
$ws.Cells.Item(6, 5).Value = '  +0.20%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.04%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.08%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.15%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.85%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.365'
$ws.Cells.Item(11, 5).Value = '  +2.21%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '3.529.96'
$ws.Cells.Item(12, 5).Value = '  +0.22%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +0.44%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +2.30%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +3.26%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Cells.Item(16, 4) '57.452.65'
$ws.Cells.Item(16, 5).Value = '  +1.18%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'Polkadot'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(17, 4) '6.21'
$ws.Cells.Item(17, 5).Value = '  +4.15%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '3.014.76'
$ws.Cells.Item(18, 5).Value = '  +0.27%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '12.83'
$ws.Cells.Item(19, 5).Value = '  +2.17%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '7.93'
$ws.Cells.Item(20, 5).Value = '  +0.62%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '328.02'
$ws.Cells.Item(21, 5).Value = '  -1.09%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '0.998'
$ws.Cells.Item(22, 5).Value = '  -0.16%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '5.67'
$ws.Cells.Item(23, 5).Value = '  -1.64%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '0.499'
$ws.Cells.Item(24, 5).Value = '  +3.15%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '64.51'
$ws.Cells.Item(25, 5).Value = '  +2.40%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -3.68%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '1.00'
$ws.Cells.Item(27, 5).Value = '  +0.07%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '0.0₃0918'
$ws.Cells.Item(28, 5).Value = '  +0.59%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '6.77'
$ws.Cells.Item(29, 5).Value = '  +0.32%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '7.33'
$ws.Cells.Item(30, 5).Value = '  +3.27%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.59%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '1.19'
$ws.Cells.Item(32, 5).Value = '  -6.32%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '20.59'
$ws.Cells.Item(33, 5).Value = '  -0.64%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +3.79%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '154.03'
$ws.Cells.Item(35, 5).Value = '  -0.09%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +3.32%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '1.27'
$ws.Cells.Item(37, 5).Value = '  -0.35%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '24.47'
$ws.Cells.Item(38, 5).Value = '  +1.19%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '0.0676'
$ws.Cells.Item(39, 5).Value = '  -0.53%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '3.047.42'
$ws.Cells.Item(40, 5).Value = '  +0.21%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '37.72'
$ws.Cells.Item(41, 5).Value = '  +1.96%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '3.85'
$ws.Cells.Item(42, 5).Value = '  +4.66%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.06%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.00%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.73%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '2.222.62'
$ws.Cells.Item(46, 5).Value = '  -2.36%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '0.980'
$ws.Cells.Item(47, 5).Value = '  -2.22%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '6.05'
$ws.Cells.Item(48, 5).Value = '  +3.57%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.15%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '19.51'
$ws.Cells.Item(50, 5).Value = '  -0.36%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '1.86'
$ws.Cells.Item(51, 5).Value = '  -6.39%  '
